$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update geocoded Latitude / Longitude values (previously placeholder 0s) ---
$ws.Range("F2").Value = 37.883536999999997
$ws.Range("G2").Value = -122.302469

$ws.Range("F3").Value = 37.865715999999999
$ws.Range("G3").Value = -122.259823

$ws.Range("F12").Value = 37.867638999999997
$ws.Range("G12").Value = -122.258082

$ws.Range("F16").Value = 37.773947
$ws.Range("G16").Value = -122.40852700000001

$ws.Range("F18").Value = 37.877206999999999
$ws.Range("G18").Value = -122.25869299999999

# --- Cells that still failed to geocode (F9/G9 stay 0) keep their "flagged" red
# fill; the rows that now have real coordinates no longer need the red
# highlight, so clear their fill back to the sheet's plain/no-fill look. Copy
# the plain format from a known unformatted cell (E19) and paste only the
# formatting so the underlying style table is reused instead of growing. ---
$ws.Range("E19").Copy() | Out-Null
$ws.Range("F2:G3").PasteSpecial(-4122)

$ws.Range("E19").Copy() | Out-Null
$ws.Range("F12:G12").PasteSpecial(-4122)

$ws.Range("E19").Copy() | Out-Null
$ws.Range("F16:G16").PasteSpecial(-4122)

$ws.Range("E19").Copy() | Out-Null
$ws.Range("F18:G18").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the view: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A16").Select() | Out-Null
$ws.Range("I7").Select() | Out-Null
